$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A22's timestamp value (microscopic correction)
$ws.Range("A22").Value = 45877.91685341435

# Append new row 23 with the new sensor reading
$ws.Range("A23").Value = 45877.95851961544
$ws.Range("B23").Value = 2025
$ws.Range("C23").Value = 32
$ws.Range("D23").Value = 15.39
$ws.Range("E23").Value = 88.17
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 4.03
$ws.Range("H23").Value = "WNW"
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = "23:00:16"

# Match the date-time number format used by column A in previous rows
$ws.Range("A23").NumberFormat = $ws.Range("A22").NumberFormat
